$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2, 3) {
    $ws.Range("D$row").Value = -0.0825
    $ws.Range("E$row").Value = -0.222

    $ws.Range("K$row").Value = 12.4
    $ws.Range("L$row").Value = 0.07745159275452843
    $ws.Range("M$row").Value = 13.6
    $ws.Range("N$row").Value = 0.04987165383204988
    $ws.Range("O$row").Value = 1.096774193548387
    $ws.Range("P$row").Value = 13.6
    $ws.Range("Q$row").Value = 0.04987165383204988
    $ws.Range("R$row").Value = 1.096774193548387

    $ws.Range("U$row").Value = 24.3
    $ws.Range("V$row").Value = 0.08910891089108912
    $ws.Range("W$row").Value = 0.03108548508398095
    $ws.Range("X$row").Value = 0.05943719655760661
    $ws.Range("Y$row").Value = -0.02835171147362566
    $ws.Range("Z$row").Value = 0.1572846055604676

    $ws.Range("AB$row").Value = 0.07230090933964295
    $ws.Range("AC$row").Value = -0.07230090933964295
    $ws.Range("AD$row").Value = 555.8
    $ws.Range("AF$row").Value = 555.8
    $ws.Range("AG$row").Value = 531.5
    $ws.Range("AH$row").Value = 0.6708509354254677
    $ws.Range("AI$row").Value = 0.6117776554760593
    $ws.Range("AJ$row").Value = 0.6609052474508829
    $ws.Range("AK$row").Value = 0.6011083465279349
}
